# Add a new bulleted list item "Audio and video files" right after the
# existing "Posting photos" item, matching its paragraph/run formatting
# (ListParagraph style, same numbered list, sz/szCs 32).

$d = $word.ActiveDocument

# Locate the last paragraph in the document ("Posting photos") and
# insert a brand-new paragraph right after it. InsertParagraphAfter()
# duplicates the originating paragraph's formatting (style, numbering,
# run properties), which is exactly what we want here.
$lastPara = $d.Paragraphs.Last
$lastPara.Range.InsertParagraphAfter()

# The document now has a new, empty trailing paragraph that inherited
# the ListParagraph/numbered-list formatting. Fill in its text.
$newPara = $d.Paragraphs.Last
$newPara.Range.Text = "Audio and video files"
